# Device-and-DeviceUse.xlsx mapping sheet update
# Commit: "Add mappings for HealthProfessional, BodyStructure, Location and Organization"
#
# The row previously holding the mapping
#   A: EHDSDeviceUse.reason   B: MedicalDevice.Indication::Problem
# (old row 26) is removed, which shifts all following rows up by one,
# bringing the HealthProfessional / Location / Organization mapping rows
# that follow it into their new positions. The trailing
# EHDSDeviceUse.reason / MedicalDevice.Indication::Diagnosis row (old row
# 29) ends up at row 28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the obsolete "MedicalDevice.Indication::Problem" mapping row;
# Excel automatically shifts subsequent rows up and drops the now-unused
# shared string on save.
$ws.Rows("26:26").Delete()

# Reflect the author's final cursor position in the saved view state.
$ws.Range("B31").Select()
